$wb = $excel.ActiveWorkbook

# Update "OFF" sheet, row 3 (R) - Week 13 logging
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 173
$wsOff.Range("C3").Value = 129
$wsOff.Range("D3").Value = 47
$wsOff.Range("E3").Value = 21
$wsOff.Range("G3").Value = 6

# Update "DEF" sheet, row 3 (R) - Week 13 logging
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 267
$wsDef.Range("C3").Value = 196
$wsDef.Range("D3").Value = 51
$wsDef.Range("E3").Value = 27
